# "Last week before rest week" - add a week of rides/workouts (rows 30-35)
# and backfill the running-total helper columns (G/H/I/J/K) on row 29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Rename the "MAF" workout titles to the new "M" naming convention ----
$ws.Range("B20").Value = "M6065 (Richmond 2015 UCI)"
$ws.Range("B21").Value = "M6070 (Volcano Loop 5)"
$ws.Range("B22").Value = "M7070 (Volcano Climb)"
$ws.Range("B23").Value = "M7570 (Knickerbocker Reverse)"
$ws.Range("B24").Value = "M7070 (Greater London Loop)"
$ws.Range("B25").Value = "M3065 rest (Innsbruckring)"
$ws.Range("B26").Value = "M8070 (Greatest London Loop)"
$ws.Range("B27").Value = "M8070 (Road to Ruins)"
$ws.Range("B28").Value = "M8070 (Kickerbocker)"

# ---- Fix up the last existing row (29): add running totals for this week ----
$ws.Range("G29").Value = 4756
$ws.Range("H29").Formula = "=I29-I28"
$ws.Range("I29").Value = 7634
$ws.Range("J29").Formula = "=K29-K28"
$ws.Range("K29").Value = 21976

# ---- Row 30: M150 (Failed Volcano Loop 10) ----
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = "2020-06-20 07:34:00"
$ws.Range("B30").Value = "M150 (Failed Volcano Loop 10)"
$ws.Range("C30").Value = "correction"
$ws.Range("D30").Value = 52.9
$ws.Range("E30").Value = 482
$ws.Range("F30").Formula = "=G30-G29"
$ws.Range("G30").Value = 4756
$ws.Range("H30").Formula = "=I30-I29"
$ws.Range("I30").Value = 7634
$ws.Range("J30").Formula = "=K30-K29"
$ws.Range("K30").Value = 23469
$ws.Range("L30").Formula = "=M30-M29"
$ws.Range("M30").Formula = "=M29"
$ws.Range("N30").Value = 160
$ws.Range("O30").Value = 83
$ws.Range("P30").Value = "02:10:19"
$ws.Range("Q30").Formula = "=TIMEVALUE(P30)*(24*60)"
$ws.Range("R30").Value = 100
$ws.Range("S30").Value = 101
$ws.Range("T30").Value = 31
$ws.Range("U30").Value = 123
$ws.Range("V30").Value = 72

# ---- Row 31: M9075 (Figure 8 Reverse) ----
$ws.Range("A31").NumberFormat = "@"
$ws.Range("A31").Value = "2020-06-22 07:34:00"
$ws.Range("B31").Value = "M9075 (Figure 8 Reverse)"
$ws.Range("C31").Value = "interval"
$ws.Range("D31").Value = 35.56
$ws.Range("E31").Value = 517
$ws.Range("F31").Value = 340
$ws.Range("G31").Value = 5096
$ws.Range("H31").Formula = "=I31-I30"
$ws.Range("I31").Value = 8168
$ws.Range("J31").Formula = "=K31-K30"
$ws.Range("K31").Value = 25145
$ws.Range("L31").Formula = "=M31-M30"
$ws.Range("M31").Value = 329170
$ws.Range("N31").Value = 160
$ws.Range("O31").Value = 83
$ws.Range("P31").Value = "1:33:55"
$ws.Range("Q31").Formula = "=TIMEVALUE(P31)*(24*60)"
$ws.Range("R31").Value = 112
$ws.Range("S31").Value = 114
$ws.Range("T31").Value = 31.8
$ws.Range("U31").Value = 126
$ws.Range("V31").Value = 65

# ---- Row 32: M9075 (Royal Pump Room 8) ----
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = "2020-06-24 07:10:00"
$ws.Range("B32").Value = "M9075 (Royal Pump Room 8)"
$ws.Range("C32").Value = "interval"
$ws.Range("D32").Value = 32
$ws.Range("E32").Value = 549
$ws.Range("F32").Value = 546
$ws.Range("G32").Value = 5642
$ws.Range("H32").Formula = "=I32-I31"
$ws.Range("I32").Value = 8721
$ws.Range("J32").Formula = "=K32-K31"
$ws.Range("K32").Value = 26825
$ws.Range("L32").Formula = "=M32-M31"
$ws.Range("M32").Value = 359031
$ws.Range("N32").Value = 160
$ws.Range("O32").Value = 82.9
$ws.Range("P32").Value = "1:40:15"
$ws.Range("Q32").Formula = "=TIMEVALUE(P32)*(24*60)"
$ws.Range("R32").Value = 112
$ws.Range("S32").Value = 114
$ws.Range("T32").Value = 31
$ws.Range("U32").Value = 124
$ws.Range("V32").Value = 63

# ---- Row 33: M9075 (Figure 8) ----
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = "2020-06-26 05:40:00"
$ws.Range("B33").Value = "M9075 (Figure 8)"
$ws.Range("C33").Value = "interval"
$ws.Range("D33").Value = 37.94
$ws.Range("E33").Value = 587
$ws.Range("F33").Formula = "=G33-G32"
$ws.Range("G33").Value = 5944
$ws.Range("H33").Formula = "=I33-I32"
$ws.Range("I33").Value = 9271
$ws.Range("J33").Formula = "=K33-K32"
$ws.Range("K33").Value = 28551
$ws.Range("L33").Formula = "=M33-M32"
$ws.Range("M33").Value = 387748
$ws.Range("N33").Value = 160
$ws.Range("O33").Value = 83.2
$ws.Range("P33").Value = "1:34:34"
$ws.Range("Q33").Formula = "=TIMEVALUE(P33)*(24*60)"
$ws.Range("R33").Value = 115
$ws.Range("S33").Value = 118
$ws.Range("T33").Value = 30
$ws.Range("U33").Value = 124
$ws.Range("V33").Value = 65

# ---- Row 34: M150 (Volcano Loop 10) ----
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value = "2020-06-27 06:39:00"
$ws.Range("B34").Value = "M150 (Volcano Loop 10)"
$ws.Range("C34").Value = "interval"
$ws.Range("D34").Value = 54.95
$ws.Range("E34").Value = 642
$ws.Range("F34").Formula = "=G34-G33"
$ws.Range("G34").Value = 6217
$ws.Range("H34").Formula = "=I34-I33"
$ws.Range("I34").Value = 9992
$ws.Range("J34").Formula = "=K34-K33"
$ws.Range("K34").Value = 30947
$ws.Range("L34").Formula = "=M34-M33"
$ws.Range("M34").Value = 426285
$ws.Range("N34").Value = 160
$ws.Range("O34").Value = 83
$ws.Range("P34").Value = "2:09:55"
$ws.Range("Q34").Formula = "=TIMEVALUE(P34)*(24*60)"
$ws.Range("R34").Value = 111
$ws.Range("S34").Value = 118
$ws.Range("T34").Value = 32
$ws.Range("U34").Value = 125
$ws.Range("V34").Value = 62

# ---- Row 35: Buy Tarmac Pro (correction entry, start of rest week) ----
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = "2020-06-27 09:10:00"
$ws.Range("B35").Value = "Buy Tarmac Pro"
$ws.Range("C35").Value = "correction"
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 642
$ws.Range("F35").Formula = "=G35-G34"
$ws.Range("G35").Value = 6217
$ws.Range("H35").Formula = "=I35-I34"
$ws.Range("I35").Value = 9992
$ws.Range("J35").Formula = "=K35-K34"
$ws.Range("K35").Value = 30947
$ws.Range("L35").Value = -408300
$ws.Range("M35").Formula = "=M34+L35"
$ws.Range("N35").Value = 160
$ws.Range("O35").Value = 83
$ws.Range("P35").NumberFormat = "@"
$ws.Range("P35").Value = "0"
$ws.Range("P35").NumberFormat = "General"
$ws.Range("Q35").Value = 0
$ws.Range("R35").Value = 0
$ws.Range("S35").Value = 0
$ws.Range("T35").Value = 0
$ws.Range("U35").Value = 0
$ws.Range("V35").Value = 0

# ---- Update the selected/active cell to reflect the new end of the sheet ----
$ws.Range("P36").Select()
